$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country list reordering ---

# India moves ahead of Paises Bajos (rows 18 and 19)
$ws.Range("A18").Value = "India"
$ws.Range("A19").Value = "Paises Bajos"

# San Vicente y las Granadinas moves ahead of Namibia (rows 194 and 195)
$ws.Range("A194").Value = "San Vicente y las Granadinas"
$ws.Range("A195").Value = "Namibia"

# --- Updated statistics ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 1170184
$ws.Range("C4").Value = 9410
$ws.Range("D4").Value = 162653
$ws.Range("E4").Value = 939529
$ws.Range("F4").Value = 16366
$ws.Range("G4").Value = 558
$ws.Range("H4").Value = 68002

# Row 7: Reino Unido
$ws.Range("B7").Value = 186599
$ws.Range("C7").Value = 4339
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 157809
$ws.Range("F7").Value = 1559
$ws.Range("G7").Value = 315
$ws.Range("H7").Value = 28446

# Row 18: India (new stats)
$ws.Range("B18").Value = 42490
$ws.Range("C18").Value = 2791
$ws.Range("D18").Value = 11775
$ws.Range("E18").Value = 29324
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 68
$ws.Range("H18").Value = 1391

# Row 19: Paises Bajos (keeps former India-row figures, unchanged from before-edit row 18)
$ws.Range("B19").Value = 40571
$ws.Range("C19").Value = 335
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 35265
$ws.Range("F19").Value = 688
$ws.Range("G19").Value = 69
$ws.Range("H19").Value = 5056
